$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = [double]"0.07324842731447992"
$ws.Range("B3").Value = [double]"-0.003795161409000828"
$ws.Range("C3").Value = [double]"0.001526093552080868"
$ws.Range("D3").Value = [double]"0.2153907562878475"
$ws.Range("E3").Value = [double]"0.07225658393248177"
$ws.Range("F3").Value = [double]"-0.006786323530935904"
$ws.Range("G3").Value = [double]"-0.0008039992870657532"
$ws.Range("H3").Value = [double]"0.06945326590547909"
$ws.Range("B4").Value = [double]"-0.0001988447026332939"
$ws.Range("C4").Value = [double]"0.003261509845689792"
$ws.Range("D4").Value = [double]"4.542246036399051"
$ws.Range("E4").Value = [double]"0.03220627872982648"
$ws.Range("F4").Value = [double]"-0.006591464608181666"
$ws.Range("G4").Value = [double]"0.006193775202915077"
$ws.Range("H4").Value = [double]"0.07304958261184663"
$ws.Range("B5").Value = [double]"0.01378693901222808"
$ws.Range("C5").Value = [double]"0.001246915074496635"
$ws.Range("D5").Value = [double]"11.129803460544"
$ws.Range("E5").Value = [double]"0.07104700122381605"
$ws.Range("F5").Value = [double]"0.01134302515284482"
$ws.Range("G5").Value = [double]"0.01623085287161133"
$ws.Range("H5").Value = [double]"0.087035366326708"
$ws.Range("B6").Value = [double]"0.02613184096272474"
$ws.Range("C6").Value = [double]"0.007241067070083312"
$ws.Range("D6").Value = [double]"11.2032524834259"
$ws.Range("E6").Value = [double]"0.1073620715344264"
$ws.Range("F6").Value = [double]"0.01193957932213652"
$ws.Range("G6").Value = [double]"0.04032410260331296"
$ws.Range("H6").Value = [double]"0.09938026827720466"
$ws.Range("B7").Value = [double]"0.1296239489556982"
$ws.Range("C7").Value = [double]"0.009305729189254507"
$ws.Range("D7").Value = [double]"13.55061898342566"
$ws.Range("E7").Value = [double]"0.1111128873428225"
$ws.Range("F7").Value = [double]"0.1113848426744075"
$ws.Range("G7").Value = [double]"0.1478630552369888"
$ws.Range("H7").Value = [double]"0.2028723762701781"
$ws.Range("B8").Value = [double]"0.06302183530241881"
$ws.Range("C8").Value = [double]"0.01099134757673806"
$ws.Range("D8").Value = [double]"9.161469941980201"
$ws.Range("E8").Value = [double]"0.107865222895561"
$ws.Range("F8").Value = [double]"0.0414791468174743"
$ws.Range("G8").Value = [double]"0.0845645237873633"
$ws.Range("H8").Value = [double]"0.1362702626168987"
$ws.Range("B9").Value = [double]"0.07659548303974745"
$ws.Range("C9").Value = [double]"0.007180895104746466"
$ws.Range("D9").Value = [double]"14.58460464864593"
$ws.Range("E9").Value = [double]"0.09498631167016788"
$ws.Range("F9").Value = [double]"0.06252115857005991"
$ws.Range("G9").Value = [double]"0.09066980750943496"
$ws.Range("H9").Value = [double]"0.1498439103542274"
$ws.Range("B10").Value = [double]"-0.07324842731447992"
$ws.Range("C10").Value = [double]"0.0004783895134771778"
$ws.Range("D10").Value = [double]"-189.0239265829394"
$ws.Range("E10").Value = [double]"1.248800245737904e-139"
$ws.Range("F10").Value = [double]"-0.07418606367300642"
$ws.Range("G10").Value = [double]"-0.07231079095595341"
$ws.Range("B11").Value = [double]"-0.02624106002961484"
$ws.Range("C11").Value = [double]"0.0005048359078661779"
$ws.Range("D11").Value = [double]"-57.52596956215765"
$ws.Range("E11").Value = [double]"0.002419921951857317"
$ws.Range("F11").Value = [double]"-0.0272305288798836"
$ws.Range("G11").Value = [double]"-0.02525159117934609"
$ws.Range("H11").Value = [double]"0.04700736728486508"
$ws.Range("B12").Value = [double]"-0.02228871391517218"
$ws.Range("C12").Value = [double]"0.0005040629494436603"
$ws.Range("D12").Value = [double]"-50.39777492183121"
$ws.Range("E12").Value = [double]"4.396786010126154e-18"
$ws.Range("F12").Value = [double]"-0.02327666854481354"
$ws.Range("G12").Value = [double]"-0.02130075928553082"
$ws.Range("H12").Value = [double]"0.05095971339930774"
$ws.Range("B13").Value = [double]"-0.01528508939017489"
$ws.Range("C13").Value = [double]"0.000489758908390348"
$ws.Range("D13").Value = [double]"-35.18080837493006"
$ws.Range("E13").Value = [double]"0.03283051483261834"
$ws.Range("F13").Value = [double]"-0.01624500775968279"
$ws.Range("G13").Value = [double]"-0.01432517102066698"
$ws.Range("H13").Value = [double]"0.05796333792430503"
$ws.Range("B14").Value = [double]"-0.01134805941195547"
$ws.Range("C14").Value = [double]"0.0004718262581276751"
$ws.Range("D14").Value = [double]"-27.55569572231164"
$ws.Range("E14").Value = [double]"6.09116909327187e-14"
$ws.Range("F14").Value = [double]"-0.0122728296759821"
$ws.Range("G14").Value = [double]"-0.01042328914792883"
$ws.Range("H14").Value = [double]"0.06190036790252446"
$ws.Range("B15").Value = [double]"-0.008123392314712445"
$ws.Range("C15").Value = [double]"0.0004663905332503993"
$ws.Range("D15").Value = [double]"-22.83870984482963"
$ws.Range("E15").Value = [double]"0.02597980898641593"
$ws.Range("F15").Value = [double]"-0.009037509060527981"
$ws.Range("G15").Value = [double]"-0.007209275568896909"
$ws.Range("H15").Value = [double]"0.06512503499976748"
$ws.Range("B16").Value = [double]"-0.007082675548082772"
$ws.Range("C16").Value = [double]"0.0004611758524982808"
$ws.Range("D16").Value = [double]"-18.04514258191906"
$ws.Range("E16").Value = [double]"0.0009354128745240235"
$ws.Range("F16").Value = [double]"-0.007986572455690937"
$ws.Range("G16").Value = [double]"-0.006178778640474608"
$ws.Range("H16").Value = [double]"0.06616575176639715"
$ws.Range("B17").Value = [double]"-0.005729584439980588"
$ws.Range("C17").Value = [double]"0.0004667201715480585"
$ws.Range("D17").Value = [double]"-14.41715546971442"
$ws.Range("E17").Value = [double]"0.05510225026191879"
$ws.Range("F17").Value = [double]"-0.00664434780424555"
$ws.Range("G17").Value = [double]"-0.004814821075715624"
$ws.Range("H17").Value = [double]"0.06751884287449933"
$ws.Range("B18").Value = [double]"-0.006160438774124211"
$ws.Range("C18").Value = [double]"0.0005067575402922984"
$ws.Range("D18").Value = [double]"-14.79965951774234"
$ws.Range("E18").Value = [double]"0.08446336132327346"
$ws.Range("F18").Value = [double]"-0.007153675689954845"
$ws.Range("G18").Value = [double]"-0.005167201858293579"
$ws.Range("H18").Value = [double]"0.06708798854035571"
$ws.Range("B19").Value = [double]"-0.004965386119004362"
$ws.Range("C19").Value = [double]"0.0004871583189600687"
$ws.Range("D19").Value = [double]"-12.39552627087503"
$ws.Range("E19").Value = [double]"0.08591827722574497"
$ws.Range("F19").Value = [double]"-0.005920208140754512"
$ws.Range("G19").Value = [double]"-0.004010564097254212"
$ws.Range("H19").Value = [double]"0.06828304119547556"
$ws.Range("B20").Value = [double]"-0.003593346238629976"
$ws.Range("C20").Value = [double]"0.0004703072846278339"
$ws.Range("D20").Value = [double]"-9.392975032980884"
$ws.Range("E20").Value = [double]"0.1067851385644758"
$ws.Range("F20").Value = [double]"-0.00451513976314539"
$ws.Range("G20").Value = [double]"-0.002671552714114561"
$ws.Range("H20").Value = [double]"0.06965508107584995"
$ws.Range("B21").Value = [double]"-0.002186109672592474"
$ws.Range("C21").Value = [double]"0.0004698457790308845"
$ws.Range("D21").Value = [double]"-5.948771536271668"
$ws.Range("E21").Value = [double]"0.007424425807992383"
$ws.Range("F21").Value = [double]"-0.003106998438004892"
$ws.Range("G21").Value = [double]"-0.001265220907180056"
$ws.Range("H21").Value = [double]"0.07106231764188745"
$ws.Range("B22").Value = [double]"-0.001480293958604121"
$ws.Range("C22").Value = [double]"0.0004650806707210647"
$ws.Range("D22").Value = [double]"-4.26244124710039"
$ws.Range("E22").Value = [double]"0.01288693306241432"
$ws.Range("F22").Value = [double]"-0.002391842927070885"
$ws.Range("G22").Value = [double]"-0.0005687449901373571"
$ws.Range("H22").Value = [double]"0.0717681333558758"
$ws.Range("B23").Value = [double]"-0.001500490375491308"
$ws.Range("C23").Value = [double]"0.0004603144120736812"
$ws.Range("D23").Value = [double]"-4.047576328841511"
$ws.Range("E23").Value = [double]"0.003841905129779102"
$ws.Range("F23").Value = [double]"-0.002402697435920454"
$ws.Range("G23").Value = [double]"-0.0005982833150621613"
$ws.Range("H23").Value = [double]"0.07174793693898861"
$ws.Range("B24").Value = [double]"-0.0009783276163937416"
$ws.Range("C24").Value = [double]"0.0004614092822246508"
$ws.Range("D24").Value = [double]"-2.500948478934665"
$ws.Range("E24").Value = [double]"0.1101709704898602"
$ws.Range("F24").Value = [double]"-0.001882681203431733"
$ws.Range("G24").Value = [double]"-7.397402935574996e-05"
$ws.Range("H24").Value = [double]"0.07227009969808618"
$ws.Range("B25").Value = [double]"-0.0006415966011542457"
$ws.Range("C25").Value = [double]"0.0004407256517024402"
$ws.Range("D25").Value = [double]"-2.074463096965729"
$ws.Range("E25").Value = [double]"0.06388096876238908"
$ws.Range("F25").Value = [double]"-0.001505410772797969"
$ws.Range("G25").Value = [double]"0.0002222175704894778"
$ws.Range("H25").Value = [double]"0.07260683071332567"
$ws.Range("B26").Value = [double]"0.07404698207497888"
$ws.Range("C26").Value = [double]"0.004399909344603928"
$ws.Range("D26").Value = [double]"22.57443277047113"
$ws.Range("E26").Value = [double]"0.002466021993036876"
$ws.Range("F26").Value = [double]"0.06542330121809975"
$ws.Range("G26").Value = [double]"0.08267066293185797"
$ws.Range("H26").Value = [double]"0.1472954093894588"
